$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''28.105.62'
$ws.Range("E2").Value = '''  -1.18%  '
$ws.Range("D3").Value = '''1.794.16'
$ws.Range("E3").Value = '''  -0.54%  '
$ws.Range("D5").Value = '''316.92'
$ws.Range("E5").Value = '''  +0.57%  '
$ws.Range("D6").Value = '''1.001'
$ws.Range("E6").Value = '''  +0.11%  '
$ws.Range("D7").Value = '''0.5361'
$ws.Range("D8").Value = '''0.3770'
$ws.Range("E8").Value = '''  -2.40%  '
$ws.Range("D9").Value = '''0.07465'
$ws.Range("E9").Value = '''  -1.86%  '
$ws.Range("D10").Value = '''41.69'
$ws.Range("E10").Value = '''  -2.05%  '
$ws.Range("E11").Value = '''  -3.02%  '
$ws.Range("E12").Value = '''  +0.12%  '
$ws.Range("E13").Value = '''  -2.90%  '
$ws.Range("D14").Value = '''6.107'
$ws.Range("E14").Value = '''  -1.38%  '
$ws.Range("B15").Value = '''WrappedEther'
$ws.Range("C15").Value = '''https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D15").Value = '''1.785.91'
$ws.Range("E15").Value = '''  -1.02%  '
$ws.Range("B16").Value = '''Chainlink'
$ws.Range("C16").Value = '''https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D16").Value = '''7.212'
$ws.Range("E16").Value = '''  -3.32%  '
$ws.Range("D17").Value = '''89.16'
$ws.Range("E17").Value = '''  -3.13%  '
$ws.Range("D18").Value = '''0.00001056'
$ws.Range("E18").Value = '''  -1.56%  '
$ws.Range("D19").Value = '''0.06461'
$ws.Range("E19").Value = '''  +0.29%  '
$ws.Range("D20").Value = '''1.000'
$ws.Range("E20").Value = '''  +0.10%  '
$ws.Range("D21").Value = '''17.31'
$ws.Range("E21").Value = '''  -0.24%  '
$ws.Range("D22").Value = '''5.902'
$ws.Range("E22").Value = '''  -1.25%  '
$ws.Range("D23").Value = '''28.131.47'
$ws.Range("E23").Value = '''  -1.10%  '
$ws.Range("E24").Value = '''  -2.17%  '
$ws.Range("D25").Value = '''2.097'
$ws.Range("E25").Value = '''  -1.66%  '
$ws.Range("D26").Value = '''154.86'
$ws.Range("E26").Value = '''  -2.64%  '
$ws.Range("D27").Value = '''20.23'
$ws.Range("E27").Value = '''  -2.39%  '
$ws.Range("D28").Value = '''1.992.92'
$ws.Range("E28").Value = '''  -0.92%  '
$ws.Range("E29").Value = '''  -5.70%  '
$ws.Range("D30").Value = '''120.68'
$ws.Range("E30").Value = '''  -2.65%  '
$ws.Range("D31").Value = '''1.116'
$ws.Range("E31").Value = '''  -0.84%  '
$ws.Range("D32").Value = '''0.1057'
$ws.Range("E32").Value = '''  +3.25%  '
$ws.Range("D33").Value = '''3.656'
$ws.Range("E33").Value = '''  -0.82%  '
$ws.Range("D34").Value = '''5.562'
$ws.Range("E34").Value = '''  -3.66%  '
$ws.Range("D35").Value = '''0.06554'
$ws.Range("E35").Value = '''  +1.56%  '
$ws.Range("D36").Value = '''0.2258'
$ws.Range("E36").Value = '''  -2.27%  '
$ws.Range("D37").Value = '''0.02284'
$ws.Range("E37").Value = '''  -1.88%  '
$ws.Range("D38").Value = '''5.020'
$ws.Range("E38").Value = '''  -3.04%  '
$ws.Range("D39").Value = '''8.439'
$ws.Range("E39").Value = '''  -4.36%  '
$ws.Range("B40").Value = '''WEMIXTOKEN'
$ws.Range("C40").Value = '''https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D40").Value = '''1.447'
$ws.Range("E40").Value = '''  +4.50%  '
$ws.Range("B41").Value = '''TheSandbox'
$ws.Range("C41").Value = '''https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D41").Value = '''0.6165'
$ws.Range("E41").Value = '''  -3.90%  '
$ws.Range("E42").Value = '''  -4.99%  '
$ws.Range("D43").Value = '''1.171'
$ws.Range("E43").Value = '''  +0.63%  '
$ws.Range("D44").Value = '''0.9998'
$ws.Range("E44").Value = '''  +0.07%  '
$ws.Range("D45").Value = '''13.37'
$ws.Range("E45").Value = '''  -1.24%  '
$ws.Range("E46").Value = '''  -0.19%  '
$ws.Range("D47").Value = '''0.5778'
$ws.Range("E47").Value = '''  -3.45%  '
$ws.Range("D48").Value = '''127.29'
$ws.Range("E48").Value = '''  +0.03%  '
$ws.Range("D49").Value = '''1.190'
$ws.Range("E49").Value = '''  +3.29%  '
$ws.Range("D50").Value = '''1.926'
$ws.Range("E50").Value = '''  -2.94%  '
$ws.Range("E51").Value = '''  -1.22%  '
